$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = "23/11/2025"
$ws.Range("B5").Value = "15:33"
$ws.Range("C5").Value = "15:33"
$ws.Range("D5").Value = "Saint-avold"
$ws.Range("E5").Value = "Jean jacque"
$ws.Range("F5").Value = "Conforme"
$ws.Range("G5").Value = "Conforme"
$ws.Range("H5").Value = "Conforme"
$ws.Range("I5").Value = "Poteau arrêt"
$ws.Range("J5").Value = "Conforme"
$ws.Range("K5").Value = "RAS"
$ws.Range("L5").Value = "casas"
$ws.Range("M5").Value = "transchool"
$ws.Range("W5").Value = "T6"
$ws.Range("X5").Value = "beau"
$ws.Range("Y5").Value = 19963
$ws.Range("Z5").Value = "Conforme"
$ws.Range("AA5").Value = "Conforme"
$ws.Range("AB5").Value = "Conforme"
$ws.Range("AC5").Value = "Conforme"
$ws.Range("AD5").Value = "Conforme"
$ws.Range("AE5").Value = "Conforme"
$ws.Range("AF5").Value = "Propre"
$ws.Range("AG5").Value = "RAS"
$ws.Range("AH5").Value = "Non observable"
$ws.Range("AI5").Value = "Non observable"
$ws.Range("AJ5").Value = "Non observable"
$ws.Range("AK5").Value = "Propre"
$ws.Range("AL5").Value = "Propre"
$ws.Range("AM5").Value = "Propre"
$ws.Range("AN5").Value = "Propre"
$ws.Range("AO5").Value = "RAS"
$ws.Range("AP5").Value = 15
$ws.Range("AQ5").Value = 3
$ws.Range("AR5").Value = "lebon"
# Row 6
$ws.Range("A6").Value = "23/11/2025"
$ws.Range("B6").Value = "15:44"
$ws.Range("C6").Value = "15:44"
$ws.Range("D6").Value = "carling"
$ws.Range("E6").Value = "LORANG"
$ws.Range("I6").Value = "Non observable"
$ws.Range("K6").Value = "ras"
$ws.Range("L6").Value = "casas"
$ws.Range("M6").Value = "transchool"
$ws.Range("W6").Value = "T12"
$ws.Range("X6").Value = "beau"
$ws.Range("Y6").Value = 25825
$ws.Range("Z6").Value = "Conforme"
$ws.Range("AA6").Value = "Conforme"
$ws.Range("AB6").Value = "Conforme"
$ws.Range("AC6").Value = "Conforme"
$ws.Range("AD6").Value = "Conforme"
$ws.Range("AE6").Value = "Conforme"
$ws.Range("AF6").Value = "Propre"
$ws.Range("AG6").Value = "ras"
$ws.Range("AH6").Value = "Conforme"
$ws.Range("AI6").Value = "Non observable"
$ws.Range("AJ6").Value = "Non observable"
$ws.Range("AK6").Value = "Propre"
$ws.Range("AL6").Value = "Propre"
$ws.Range("AM6").Value = "Propre"
$ws.Range("AN6").Value = "Propre"
$ws.Range("AO6").Value = "ras"
$ws.Range("AP6").Value = 20
$ws.Range("AQ6").Value = 3
$ws.Range("AR6").Value = "lebon"
# Row 7
$ws.Range("A7").Value = "24/11/2025"
$ws.Range("B7").Value = "06:49"
$ws.Range("C7").Value = "06:01"
$ws.Range("D7").Value = "Merlebach"
$ws.Range("E7").Value = "Friderich"
$ws.Range("I7").Value = "Non observable"
$ws.Range("K7").Value = "ras"
$ws.Range("L7").Value = "rgeFluo57"
$ws.Range("N7").Value = "Sa"
$ws.Range("T7").Value = "SA24"
$ws.Range("X7").Value = "beau"
$ws.Range("Y7").Value = 22320
$ws.Range("Z7").Value = "Conforme"
$ws.Range("AA7").Value = "Conforme"
$ws.Range("AB7").Value = "Conforme"
$ws.Range("AC7").Value = "Conforme"
$ws.Range("AD7").Value = "Conforme"
$ws.Range("AE7").Value = "Conforme"
$ws.Range("AF7").Value = "Propre"
$ws.Range("AG7").Value = "ras"
$ws.Range("AH7").Value = "Non observable"
$ws.Range("AI7").Value = "Non observable"
$ws.Range("AJ7").Value = "Conforme"
$ws.Range("AK7").Value = "Propre"
$ws.Range("AL7").Value = "Propre"
$ws.Range("AM7").Value = "Propre"
$ws.Range("AN7").Value = "Propre"
$ws.Range("AO7").Value = "ras"
$ws.Range("AP7").Value = 10
$ws.Range("AQ7").Value = 2
$ws.Range("AR7").Value = "lebon"
